$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Mdk","Sdc1","ECs",3,1,2.512729,7.538187,0.02190726325199687,0.02190726325199687,3,1,0.9159313333333333,2.747794,0.03641350786393945,0.03641350786393944,2.301487223275333,20.713385009478,0.0007977203027039797,0.0007977203027039795),
  @("ECs","Mdk","Sdc1","FAPs",3,1,2.512729,7.538187,0.02190726325199687,0.02190726325199687,3,1,10.56834633333333,31.705039,0.4201521973455823,0.4201521973455822,26.55539031381033,238.998512824293,0.009204384793154611,0.009204384793154609),
  @("ECs","Mdk","Sdc1","M2",3,1,2.512729,7.538187,0.02190726325199687,0.02190726325199687,2,0.6666666666666666,1.220967333333333,3.662902,0.04854043308262539,0.04854043308262537,3.067960026519333,27.611640238674,0.001063388045907012,0.001063388045907012),
  @("ECs","Mdk","Sdc1","sCs",3,1,2.512729,7.538187,0.02190726325199687,0.02190726325199687,3,1,12.44836933333333,37.345108,0.494893861707853,0.4948938617078529,31.27937862657733,281.514407639196,0.01084177011023127,0.01084177011023127),
  @("FAPs","Mdk","Sdc1","ECs",3,1,92.89399466666667,278.681984,0.8098976036382196,0.8098976036382197,3,1,0.9159313333333333,2.747794,0.03641350786393945,0.03641350786393944,85.08452039369956,765.760683543296,0.02949121275906603,0.02949121275906602),
  @("FAPs","Mdk","Sdc1","FAPs",3,1,92.89399466666667,278.681984,0.8098976036382196,0.8098976036382197,3,1,10.56834633333333,31.705039,0.4201521973455823,0.4201521973455822,981.7359079241529,8835.623171317377,0.3402802577935194,0.3402802577935194),
  @("FAPs","Mdk","Sdc1","M2",3,1,92.89399466666667,278.681984,0.8098976036382196,0.8098976036382197,2,0.6666666666666666,1.220967333333333,3.662902,0.04854043308262539,0.04854043308262537,113.4205329508409,1020.784796557568,0.03931278043317966,0.03931278043317965),
  @("FAPs","Mdk","Sdc1","sCs",3,1,92.89399466666667,278.681984,0.8098976036382196,0.8098976036382197,3,1,12.44836933333333,37.345108,0.494893861707853,0.4948938617078529,1156.378754459364,10407.40879013427,0.4008133526524546,0.4008133526524545),
  @("M2","Mdk","Sdc1","ECs",3,1,1.610639333333333,4.831918,0.0140423817607685,0.0140423817607685,3,1,0.9159313333333333,2.747794,0.03641350786393945,0.03641350786393944,1.475235032099111,13.277115288892,0.0005113323786741837,0.0005113323786741836),
  @("M2","Mdk","Sdc1","FAPs",3,1,1.610639333333333,4.831918,0.0140423817607685,0.0140423817607685,3,1,10.56834633333333,31.705039,0.4201521973455823,0.4201521973455822,17.02179429275578,153.196148634802,0.005899937552752412,0.005899937552752412),
  @("M2","Mdk","Sdc1","M2",3,1,1.610639333333333,4.831918,0.0140423817607685,0.0140423817607685,2,0.6666666666666666,1.220967333333333,3.662902,0.04854043308262539,0.04854043308262537,1.966538011781778,17.698842106036,0.0006816232921792626,0.0006816232921792625),
  @("M2","Mdk","Sdc1","sCs",3,1,1.610639333333333,4.831918,0.0140423817607685,0.0140423817607685,3,1,12.44836933333333,37.345108,0.494893861707853,0.4948938617078529,20.04983328412711,180.448499557144,0.006949488537162644,0.006949488537162643),
  @("sCs","Mdk","Sdc1","ECs",3,1,17.68108066666667,53.04324200000001,0.154152751349015,0.154152751349015,3,1,0.9159313333333333,2.747794,0.03641350786393945,0.03641350786393944,16.19465578979423,145.751902108148,0.005613242423495261,0.005613242423495258),
  @("sCs","Mdk","Sdc1","FAPs",3,1,17.68108066666667,53.04324200000001,0.154152751349015,0.154152751349015,3,1,10.56834633333333,31.705039,0.4201521973455823,0.4201521973455822,186.8597840329376,1681.738056296438,0.06476761720615583,0.0647676172061558),
  @("sCs","Mdk","Sdc1","M2",3,1,17.68108066666667,53.04324200000001,0.154152751349015,0.154152751349015,2,0.6666666666666666,1.220967333333333,3.662902,0.04854043308262539,0.04854043308262537,21.58802191203156,194.292197208284,0.007482641311359452,0.00748264131135945),
  @("sCs","Mdk","Sdc1","sCs",3,1,17.68108066666667,53.04324200000001,0.154152751349015,0.154152751349015,3,1,12.44836933333333,37.345108,0.494893861707853,0.4948938617078529,220.1006223511263,1980.905601160136,0.07628925040800447,0.07628925040800445)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}

Write-Output "done"